$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "H2-K1"
$ws.Range("C2").Value = "Erbb2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 282.5164233333333
$ws.Range("H2").Value = 847.54927
$ws.Range("I2").Value = 0.7504954445259187
$ws.Range("J2").Value = 0.7504954445259185
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 1.534538333333333
$ws.Range("N2").Value = 4.603615
$ws.Range("O2").Value = 0.1494637976135089
$ws.Range("P2").Value = 0.1494637976135089
$ws.Range("Q2").Value = 433.5322814012277
$ws.Range("R2").Value = 3901.790532611049
$ws.Range("S2").Value = 0.1121718992304823
$ws.Range("T2").Value = 0.1121718992304823

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "H2-K1"
$ws.Range("C3").Value = "Erbb2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 282.5164233333333
$ws.Range("H3").Value = 847.54927
$ws.Range("I3").Value = 0.7504954445259187
$ws.Range("J3").Value = 0.7504954445259185
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.973328333333334
$ws.Range("N3").Value = 14.919985
$ws.Range("O3").Value = 0.4844014146353658
$ws.Range("P3").Value = 0.4844014146353658
$ws.Range("Q3").Value = 1405.046932795661
$ws.Range("R3").Value = 12645.42239516095
$ws.Range("S3").Value = 0.3635410550057527
$ws.Range("T3").Value = 0.3635410550057526

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "H2-K1"
$ws.Range("C4").Value = "Erbb2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 282.5164233333333
$ws.Range("H4").Value = 847.54927
$ws.Range("I4").Value = 0.7504954445259187
$ws.Range("J4").Value = 0.7504954445259185
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.75909
$ws.Range("N4").Value = 11.27727
$ws.Range("O4").Value = 0.3661347877511252
$ws.Range("P4").Value = 0.3661347877511252
$ws.Range("Q4").Value = 1062.0046617881
$ws.Range("R4").Value = 9558.0419560929
$ws.Range("S4").Value = 0.2747824902896835
$ws.Range("T4").Value = 0.2747824902896835

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "H2-K1"
$ws.Range("C5").Value = "Erbb2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 68.514867
$ws.Range("H5").Value = 205.544601
$ws.Range("I5").Value = 0.1820074562714184
$ws.Range("J5").Value = 0.1820074562714184
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 1.534538333333333
$ws.Range("N5").Value = 4.603615
$ws.Range("O5").Value = 0.1494637976135089
$ws.Range("P5").Value = 0.1494637976135089
$ws.Range("Q5").Value = 105.138689814735
$ws.Range("R5").Value = 946.2482083326149
$ws.Range("S5").Value = 0.02720352560830086
$ws.Range("T5").Value = 0.02720352560830086

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "H2-K1"
$ws.Range("C6").Value = "Erbb2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 68.514867
$ws.Range("H6").Value = 205.544601
$ws.Range("I6").Value = 0.1820074562714184
$ws.Range("J6").Value = 0.1820074562714184
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 4.973328333333334
$ws.Range("N6").Value = 14.919985
$ws.Range("O6").Value = 0.4844014146353658
$ws.Range("P6").Value = 0.4844014146353658
$ws.Range("Q6").Value = 340.746929305665
$ws.Range("R6").Value = 3066.722363750985
$ws.Range("S6").Value = 0.08816466929205956
$ws.Range("T6").Value = 0.08816466929205954

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "H2-K1"
$ws.Range("C7").Value = "Erbb2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 68.514867
$ws.Range("H7").Value = 205.544601
$ws.Range("I7").Value = 0.1820074562714184
$ws.Range("J7").Value = 0.1820074562714184
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.75909
$ws.Range("N7").Value = 11.27727
$ws.Range("O7").Value = 0.3661347877511252
$ws.Range("P7").Value = 0.3661347877511252
$ws.Range("Q7").Value = 257.55355139103
$ws.Range("R7").Value = 2317.98196251927
$ws.Range("S7").Value = 0.06663926137105797
$ws.Range("T7").Value = 0.06663926137105797

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "H2-K1"
$ws.Range("C8").Value = "Erbb2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 25.40860066666666
$ws.Range("H8").Value = 76.22580199999999
$ws.Range("I8").Value = 0.06749709920266306
$ws.Range("J8").Value = 0.06749709920266304
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 1.534538333333333
$ws.Range("N8").Value = 4.603615
$ws.Range("O8").Value = 0.1494637976135089
$ws.Range("P8").Value = 0.1494637976135089
$ws.Range("Q8").Value = 38.99047171935887
$ws.Range("R8").Value = 350.9142454742299
$ws.Range("S8").Value = 0.01008837277472577
$ws.Range("T8").Value = 0.01008837277472577

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "H2-K1"
$ws.Range("C9").Value = "Erbb2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 25.40860066666666
$ws.Range("H9").Value = 76.22580199999999
$ws.Range("I9").Value = 0.06749709920266306
$ws.Range("J9").Value = 0.06749709920266304
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 4.973328333333334
$ws.Range("N9").Value = 14.919985
$ws.Range("O9").Value = 0.4844014146353658
$ws.Range("P9").Value = 0.4844014146353658
$ws.Range("Q9").Value = 126.3653136058855
$ws.Range("R9").Value = 1137.28782245297
$ws.Range("S9").Value = 0.03269569033755361
$ws.Range("T9").Value = 0.0326956903375536

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "H2-K1"
$ws.Range("C10").Value = "Erbb2"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 25.40860066666666
$ws.Range("H10").Value = 76.22580199999999
$ws.Range("I10").Value = 0.06749709920266306
$ws.Range("J10").Value = 0.06749709920266304
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.75909
$ws.Range("N10").Value = 11.27727
$ws.Range("O10").Value = 0.3661347877511252
$ws.Range("P10").Value = 0.3661347877511252
$ws.Range("Q10").Value = 95.51321668005998
$ws.Range("R10").Value = 859.6189501205398
$ws.Range("S10").Value = 0.02471303609038368
$ws.Range("T10").Value = 0.02471303609038367

